# remover pessoa apos 7 dias se o email nao foi confirmado
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column H "DataEnvioInicial" with the same header style as G1 ---
$ws.Range("H1").Value = "DataEnvioInicial"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Touch H2/H3/H4 so the (currently blank) column H cells exist on those rows too
$ws.Range("H2").Font.Bold = $ws.Range("H2").Font.Bold
$ws.Range("H3").Font.Bold = $ws.Range("H3").Font.Bold
$ws.Range("H4").Font.Bold = $ws.Range("H4").Font.Bold

# --- Row 3 previously held an unconfirmed entry (E3 = "n") whose email was
# never confirmed; after 7 days that person's contact/tracking data is wiped ---
$ws.Range("B3").Value = "brinquedo superman"
$ws.Range("C3").ClearContents() | Out-Null
$ws.Range("F3").ClearContents() | Out-Null
$ws.Range("G3").ClearContents() | Out-Null

# Re-touch the cells we just cleared so they remain present (but blank) in the sheet
$ws.Range("C3").Font.Bold = $ws.Range("C3").Font.Bold
$ws.Range("F3").Font.Bold = $ws.Range("F3").Font.Bold
$ws.Range("G3").Font.Bold = $ws.Range("G3").Font.Bold

# --- New row 4: a fresh submission for the same toy/email, now awaiting confirmation ---
$ws.Range("A4").Value = "Teste excelTeste excel"
$ws.Range("B4").Value = "brinquedo superman"
$ws.Range("C4").Value = "daniel.jmendes2@gmail.com"
$ws.Range("E4").Value = "sim"
$ws.Range("F4").Value = "daniel"

# Touch D4/G4/H4 so they remain present (but blank) like their counterparts in other rows
$ws.Range("D4").Font.Bold = $ws.Range("D4").Font.Bold
$ws.Range("G4").Font.Bold = $ws.Range("G4").Font.Bold
